# The generated colo-location table had a stale row ("TSN" / Tianjin, China)
# sitting ahead of its correct alphabetical position. The refreshed export
# drops that row entirely, so every row below it (NRT onward) shifts up by
# one — ending with the table one row shorter (A1:H333 instead of A1:H334).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 235 (colo "TSN", Tianjin, China) and shift everything
# below it up by one row, matching the regenerated data dump.
$ws.Rows(235).Delete()
